# "finish Calendar & swap by function"
#
# - Typography sheet: extend the Wildcard Ranges for Typography_00 with 0-9,
#   and add "%" as a Wildcard Character for Typography_01.
# - Translation sheet: rename the GB text used by SingleUseId16 from TEMP to
#   YEAR (the old TEMP text moves to a brand-new SingleUseId row), and add
#   five brand-new rows finishing out the Calendar / Humidity pieces
#   (rows 140-144).

$wb = $excel.ActiveWorkbook

# --- Typography sheet ---
$typo = $wb.Worksheets.Item("Typography")

# Row 7 (Typography_00): Wildcard Ranges changes from "a-z,A-Z" to "a-z,A-Z,0-9"
$typo.Range("H7").Value = "a-z,A-Z,0-9"

# Row 8 (Typography_01): Wildcard Characters changes from empty to "%"
$typo.Range("G8").Value = "%"

# --- Translation sheet ---
$trans = $wb.Worksheets.Item("Translation")

# Row 17 (SingleUseId16): GB text changes from "TEMP" to "YEAR"
$trans.Range("E17").Value = "YEAR"

# New row 140: SingleUseId139 / Large / Left / TEMP / LTR
$trans.Range("B140").Value = "SingleUseId139"
$trans.Range("C140").Value = "Large"
$trans.Range("D140").Value = "Left"
$trans.Range("E140").Value = "TEMP"
$trans.Range("F140").Value = "LTR"
$trans.Range("B140:F140").Style = "Normal"

# New row 141: SingleUseId140 / Large / Left / HUMID / LTR
$trans.Range("B141").Value = "SingleUseId140"
$trans.Range("C141").Value = "Large"
$trans.Range("D141").Value = "Left"
$trans.Range("E141").Value = "HUMID"
$trans.Range("F141").Value = "LTR"
$trans.Range("B141:F141").Style = "Normal"

# New row 142: SingleUseId141 / Typography_01 / Right / <v> / LTR
$trans.Range("B142").Value = "SingleUseId141"
$trans.Range("C142").Value = "Typography_01"
$trans.Range("D142").Value = "Right"
$trans.Range("E142").Value = "<v>"
$trans.Range("F142").Value = "LTR"
$trans.Range("B142:F142").Style = "Normal"

# New row 143: SingleUseId143 / Large / Left / % / LTR
$trans.Range("B143").Value = "SingleUseId143"
$trans.Range("C143").Value = "Large"
$trans.Range("D143").Value = "Left"
$trans.Range("E143").Value = "%"
$trans.Range("F143").Value = "LTR"
$trans.Range("B143:F143").Style = "Normal"

# New row 144: SingleUseId144 / Typography_01 / Left / 0 / LTR
# ("0" must stay text, matching the other rows that use the literal "0" text)
$trans.Range("B144").Value = "SingleUseId144"
$trans.Range("C144").Value = "Typography_01"
$trans.Range("D144").Value = "Left"
$trans.Range("E144").NumberFormat = "@"
$trans.Range("E144").Value = "0"
$trans.Range("F144").Value = "LTR"
$trans.Range("B144:F144").Style = "Normal"
